$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the four new exchange-rate header columns (K1:N1)
$ws.Range("K1").Value = "From Currency"
$ws.Range("L1").Value = "To Currency"
$ws.Range("M1").Value = "Exchange Rate "
$ws.Range("N1").Value = "As Of"

# The new header cells pick up a distinct (but visually "Normal") style -
# nudge formatting so the engine allocates a new cellXf for them, matching
# the workbook's updated style table.
$ws.Range("K1:N1").WrapText = $false

# Move / collapse the active selection to K2, as in the edited workbook
$ws.Range("K2").Select()
